$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "42.856.72") are not coerced into numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "42.856.72"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "2.215.58"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "254.29"
$ws.Range("E5").Value = "  +3.41%  "

# Row 6
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  -0.53%  "

# Row 7
$ws.Range("D7").Value = "75.49"
$ws.Range("E7").Value = "  -0.60%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -4.46%  "

# Row 10
$ws.Range("D10").Value = "41.18"
$ws.Range("E10").Value = "  +0.23%  "

# Row 11
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  -2.33%  "

# Row 12
$ws.Range("D12").Value = "6.88"
$ws.Range("E12").Value = "  -1.37%  "

# Row 13
$ws.Range("E13").Value = "  +0.47%  "

# Row 14
$ws.Range("D14").Value = "2.544.76"
$ws.Range("E14").Value = "  -0.53%  "

# Row 15
$ws.Range("D15").Value = "14.35"
$ws.Range("E15").Value = "  -1.99%  "

# Row 16
$ws.Range("D16").Value = "2.216.82"
$ws.Range("E16").Value = "  -1.25%  "

# Row 17
$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  -3.29%  "

# Row 18
$ws.Range("D18").Value = "42.752.14"
$ws.Range("E18").Value = "  -0.22%  "

# Row 19
$ws.Range("E19").Value = "  -2.75%  "

# Row 20
$ws.Range("D20").Value = "71.01"
$ws.Range("E20").Value = "  -0.46%  "

# Row 21
$ws.Range("D21").Value = "5.94"
$ws.Range("E21").Value = "  -1.11%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "229.29"
$ws.Range("E22").Value = "  -0.89%  "

# Row 23
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "2.18"
$ws.Range("E23").Value = "  -2.14%  "

# Row 24
$ws.Range("D24").Value = "9.31"
$ws.Range("E24").Value = "  -8.65%  "

# Row 25
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("D26").Value = "10.57"
$ws.Range("E26").Value = "  -2.86%  "

# Row 27
$ws.Range("D27").Value = "3.38"
$ws.Range("E27").Value = "  -1.21%  "

# Row 28
$ws.Range("D28").Value = "38.50"
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.71%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  -2.91%  "

# Row 31
$ws.Range("D31").Value = "172.76"
$ws.Range("E31").Value = "  -0.24%  "

# Row 32
$ws.Range("D32").Value = "20.23"
$ws.Range("E32").Value = "  -0.44%  "

# Row 33
$ws.Range("D33").Value = "0.0844"
$ws.Range("E33").Value = "  +6.64%  "

# Row 34
$ws.Range("D34").Value = "5.20"
$ws.Range("E34").Value = "  -2.74%  "

# Row 35
$ws.Range("E35").Value = "  -1.34%  "

# Row 36
$ws.Range("E36").Value = "  -1.67%  "

# Row 37
$ws.Range("D37").Value = "0.0347"
$ws.Range("E37").Value = "  +6.13%  "

# Row 38
$ws.Range("D38").Value = "4.28"
$ws.Range("E38").Value = "  -1.75%  "

# Row 39
$ws.Range("D39").Value = "12.38"
$ws.Range("E39").Value = "  -3.67%  "

# Row 40
$ws.Range("D40").Value = "2.10"
$ws.Range("E40").Value = "  -2.02%  "

# Row 41
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  +17.52%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "59.99"
$ws.Range("E42").Value = "  +0.22%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.197"
$ws.Range("E43").Value = "  -3.36%  "

# Row 44
$ws.Range("D44").Value = "5.25"
$ws.Range("E44").Value = "  -5.83%  "

# Row 45
$ws.Range("D45").Value = "101.73"
$ws.Range("E45").Value = "  -4.77%  "

# Row 46
$ws.Range("D46").Value = "8.33"
$ws.Range("E46").Value = "  -3.91%  "

# Row 47
$ws.Range("D47").Value = "0.0976"
$ws.Range("E47").Value = "  -1.24%  "

# Row 48
$ws.Range("D48").Value = "0.458"
$ws.Range("E48").Value = "  +2.32%  "

# Row 49
$ws.Range("E49").Value = "  -0.22%  "

# Row 50
$ws.Range("E50").Value = "  -1.10%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.438.12"
$ws.Range("E51").Value = "  -0.31%  "
